$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain-text column updates (Coin / Link / Volume) ---
# Row 18's "worst in 24h" suffix was dropped, and rows 41-43 got rotated
# (KickToken -> BKEXToken -> CEJI -> KickToken). None of these values look
# like numbers, so direct assignment keeps them as text, exactly as before.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Price (col D) and Hora (col G) updates ---
# These values look numeric, so Excel COM would silently coerce them to
# doubles (losing trailing zeros / exact text) unless the range is
# pre-formatted as Text. Set NumberFormat="@" first, write the literal
# string, then restore the "Normal" style so the stored cell keeps the
# workbook's original (unstyled) appearance.
$numericCells = @{
    "D2" = "243.73"
    "G2" = "10"
    "D3" = "23.77"
    "G3" = "10"
    "D4" = "5.252"
    "G4" = "10"
    "D5" = "0.05812"
    "G5" = "10"
    "D6" = "6.479"
    "G6" = "10"
    "D7" = "3.354"
    "G7" = "10"
    "D8" = "0.8082"
    "G8" = "10"
    "D9" = "0.8782"
    "G9" = "10"
    "D10" = "0.1388"
    "G10" = "10"
    "D11" = "0.07263"
    "G11" = "10"
    "D12" = "0.03075"
    "G12" = "10"
    "G13" = "10"
    "D14" = "0.09315"
    "G14" = "10"
    "D15" = "3.856"
    "G15" = "10"
    "G16" = "10"
    "D17" = "0.04687"
    "G17" = "10"
    "D18" = "0.0006018"
    "G18" = "10"
    "D19" = "0.006182"
    "G19" = "10"
    "D20" = "0.001265"
    "G20" = "10"
    "D21" = "0.004593"
    "G21" = "10"
    "G22" = "10"
    "G23" = "10"
    "D24" = "2.158"
    "G24" = "10"
    "D25" = "0.3209"
    "G25" = "10"
    "D26" = "0.1310"
    "G26" = "10"
    "G27" = "10"
    "G28" = "10"
    "G29" = "10"
    "G30" = "10"
    "G31" = "10"
    "G32" = "10"
    "G33" = "10"
    "G34" = "10"
    "G35" = "10"
    "G36" = "10"
    "G37" = "10"
    "G38" = "10"
    "G39" = "10"
    "D40" = "0.03795"
    "G40" = "10"
    "D41" = "0.1054"
    "G41" = "10"
    "D42" = "0.002439"
    "G42" = "10"
    "D43" = "0.003230"
    "G43" = "10"
    "D44" = "0.007980"
    "G44" = "10"
    "D45" = "0.00005496"
    "G45" = "10"
    "G46" = "10"
    "D47" = "0.5199"
    "G47" = "10"
    "D48" = "0.007019"
    "G48" = "10"
    "D49" = "0.00002100"
    "G49" = "10"
    "D50" = "0.0002000"
    "G50" = "10"
    "G51" = "10"
}
foreach ($addr in $numericCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericCells[$addr]
    $cell.Style = "Normal"
}
